$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared / rich-text string edits ---
# A8: "Volume 31   Number  46" -> "...47" (report/volume number bump)
$ws.Range("A8").Characters(21, 2).Text = "47"

# C9: "Report Covering the Week  11/11/2024  Through  11/17/2024" -> next week
$ws.Range("C9").Characters(27, 10).Text = "11/18/2024"
$ws.Range("C9").Characters(48, 10).Text = "11/24/2024"

# --- Row 18 / Row 20: C/D pair becomes a "no prior data" placeholder ---
# (numeric 0 divisor cases rendered as shared text "0" / "***.*", General-format
#  style matching the existing placeholder rows, e.g. row 22)
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "'***.*"
$ws.Range("D20").Value = "'0"
$ws.Range("E20").Value = "'***.*"
$ws.Range("D22:E22").Copy() | Out-Null
$ws.Range("D18:E18").PasteSpecial(-4122) | Out-Null
$ws.Range("D20:E20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Weekly crime-stat figures (new data collected) ---
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 23
$ws.Range("K16").Value = 27.777777777777
$ws.Range("L16").Value = 27.777777777777
$ws.Range("M16").Value = 4.545454545454
$ws.Range("N16").Value = -61.666666666666
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = -25
$ws.Range("J17").Value = 82
$ws.Range("K17").Value = 2.439024390243
$ws.Range("M17").Value = 90.909090909090
$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 200
$ws.Range("I18").Value = 40
$ws.Range("K18").Value = -18.367346938775
$ws.Range("L18").Value = 8.108108108108
$ws.Range("M18").Value = -58.333333333333
$ws.Range("N18").Value = -86.970684039087
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -25
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -57.142857142857
$ws.Range("I19").Value = 240
$ws.Range("J19").Value = 264
$ws.Range("K19").Value = -9.090909090909
$ws.Range("L19").Value = -2.040816326530
$ws.Range("M19").Value = 79.104477611940
$ws.Range("N19").Value = 8.597285067873
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 36
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -68.141592920354
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -94.435857805255
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = 29
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = -34.090909090909
$ws.Range("I21").Value = 431
$ws.Range("J21").Value = 487
$ws.Range("K21").Value = -11.498973305954
$ws.Range("L21").Value = -7.311827956989
$ws.Range("M21").Value = 28.273809523809
$ws.Range("N21").Value = -67.859806114839
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 200
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = -14.893617021276
$ws.Range("I24").Value = 340
$ws.Range("J24").Value = 444
$ws.Range("K24").Value = -23.423423423423
$ws.Range("L24").Value = -23.076923076923
$ws.Range("M24").Value = -31.863727454909
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 25
$ws.Range("H25").Value = 177.777777777778
$ws.Range("I25").Value = 158
$ws.Range("J25").Value = 187
$ws.Range("K25").Value = -15.508021390374
$ws.Range("L25").Value = 47.663551401869
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 80
$ws.Range("F26").Value = 17
$ws.Range("H26").Value = 88.888888888888
$ws.Range("I26").Value = 150
$ws.Range("J26").Value = 175
$ws.Range("K26").Value = -14.285714285714
$ws.Range("L26").Value = -11.242603550295
$ws.Range("M26").Value = -23.076923076923

Write-Host "Applied weekly crime data update."
